# Update cohort demographics figures: add match analysis for Males and Females.
# Each Find/Replace targets a unique text run in the document, so ordering is
# not significant.

$d = $word.ActiveDocument

$replacements = @(
    @("N = 1,996", "N = 1,995"),
    @(", N = 1206", ", N = 1205"),
    @("1,182 (59.22%)", "1,181 (59.20%)"),
    @("717 (59.45%)", "716 (59.42%)"),
    @("814 (40.78%)", "814 (40.80%)"),
    @("489 (40.55%)", "489 (40.58%)"),
    @("1,219 (61.07%)", "1,219 (61.10%)"),
    @("913 (75.70%)", "913 (75.77%)"),
    @("777 (38.93%)", "776 (38.90%)"),
    @("293 (24.30%)", "292 (24.23%)"),
    @("981 (49.15%)", "981 (49.17%)"),
    @("243 (20.15%)", "243 (20.17%)"),
    @("1,637 (82.01%)", "1,636 (82.01%)"),
    @("1,186 (98.34%)", "1,185 (98.34%)"),
    @("3.25 (1.56)", "3.26 (1.56)")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
